# Clean repo pass 1
#
# 1. Fill in the missing "Species" (column C) values on the T.scriptus sheet
#    with "Tragelaphus scriptus" for every data row that doesn't already
#    have it.
# 2. Fix the C.spadix!_FilterDatabase defined name so the sheet name is
#    quoted (it contains a "." which Excel re-quotes on save).
# 3. Re-point the active sheet / selection to T.scriptus (matches the new
#    activeTab + tabSelected + selection in the diff).

$wb = $excel.ActiveWorkbook

# --- 1. Backfill the Species column on T.scriptus -------------------------
$ws6 = $wb.Worksheets.Item("T.scriptus")

$alreadyFilled = @(15, 16, 22, 46, 65, 66, 72, 76, 82, 105)

for ($r = 2; $r -le 121; $r++) {
    if ($alreadyFilled -contains $r) { continue }
    $ws6.Cells.Item($r, 3).Value = "Tragelaphus scriptus"
}

# --- 2. Fix the C.spadix defined name quoting ------------------------------
foreach ($n in $wb.Names) {
    if ($n.Name -eq "C.spadix!_FilterDatabase") {
        $n.RefersTo = "='C.spadix'!`$A`$1:`$J`$121"
    }
}

# --- 3. Make T.scriptus the active/selected sheet --------------------------
$ws6.Activate()
$ws6.Range("Q16").Select()
